$d = $word.ActiveDocument

# Shared run-properties block used by all the runs in these two paragraphs.
$rPr = '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr>'

# Shared paragraph-properties block (bulleted list item) used by both paragraphs.
$pPr = '<w:pPr><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="32"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="468"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr>'

$ns = ' xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Paragraph 1: "Issues raised by Customers via Exor Support"
# -> drop the two proofErr markers, rename the middle run's text to "Bentley".
$p1Xml = '<w:p' + $ns + ' w:rsidR="00853BDA" w:rsidRDefault="00C64E4F">' + $pPr + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve">Issues raised by Customers via </w:t></w:r>' + `
    '<w:r>' + $rPr + '<w:t>Bentley</w:t></w:r>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> Support</w:t></w:r>' + `
    '</w:p>'

# Paragraph 2: "Issues raised internally by Exor"
# -> drop the two proofErr markers, rename the final run's text to "Bentley".
$p2Xml = '<w:p' + $ns + ' w:rsidR="00853BDA" w:rsidRDefault="00C64E4F">' + $pPr + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve">Issues raised internally by </w:t></w:r>' + `
    '<w:r>' + $rPr + '<w:t>Bentley</w:t></w:r>' + `
    '</w:p>'

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Contains("Issues raised by Customers via Exor")) {
        $p.Range.InsertXML($p1Xml)
        break
    }
}

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Contains("Issues raised internally by Exor")) {
        $p.Range.InsertXML($p2Xml)
        break
    }
}
